# Consumption_Forecast_Historical.xlsx — roll the forecast window forward by one day.
#
# Layout: row 1 = header; rows 2-97 = 96 quarter-hour slots for "day 1"
# (25.01.2026); rows 98-193 = 96 quarter-hour slots for "day 2" (26.01.2026).
# Columns: A Timestamp, B Forecasted Consumption (MW), C Quarter (1..96,
# repeats per day — unaffected by the roll), D Lookup ("DD.MM.YYYY" + quarter
# index, used as a helper key).
#
# The retrained model drops the oldest day and appends a freshly forecast
# day: old day-2 (rows 98-193) slides up to become the new day-1 (rows
# 2-97), and a brand-new day-3 (27.01.2026) is written into rows 98-193 with
# updated forecast figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 2
$dayRows = 96
$secondDayStart = $firstDataRow + $dayRows      # 98
$lastDataRow = $secondDayStart + $dayRows - 1   # 193

# --- 1) Snapshot the current "day 2" block (rows 98-193) before overwriting anything ---
$oldA = @()
$oldB = @()
$oldD = @()
for ($i = 0; $i -lt $dayRows; $i++) {
    $r = $secondDayStart + $i
    $oldA += $ws.Cells.Item($r, 1).Value()
    $oldB += $ws.Cells.Item($r, 2).Value()
    $oldD += $ws.Cells.Item($r, 4).Value()
}

# --- 2) Move that block up into rows 2-97 (becomes the new "day 1") ---
# Column C (Quarter index 1..96) is identical for every day, so it is left untouched.
for ($i = 0; $i -lt $dayRows; $i++) {
    $r = $firstDataRow + $i
    $ws.Cells.Item($r, 1).Value = $oldA[$i]
    $ws.Cells.Item($r, 2).Value = $oldB[$i]
    $ws.Cells.Item($r, 4).Value = $oldD[$i]
}

# --- 3) Freshly forecast new "day 3" (27.01.2026) values for rows 98-193 ---
$newConsumption = @(6120,6090,6060,6030,5990,0,5980,0,5990,6000,6010,6020,6030,6040,6050,6090,6160,6250,6350,6490,6660,6840,7030,7230,7430,7620,7810,7980,8140,8270,8380,8470,8530,8580,8600,0,8590,8570,8550,8530,8510,8490,8480,8460,8450,8440,0,0,8450,0,0,8440,8400,8380,8370,8360,8340,8320,0,0,0,8330,8350,8370,8390,8410,8440,8480,8530,8550,0,0,8500,8460,8430,8400,8340,8270,8220,8130,8000,7870,7770,7640,7470,7320,7160,7020,6870,6720,6630,6520,6580,6520,6460,6420)

for ($i = 0; $i -lt $dayRows; $i++) {
    $r = $secondDayStart + $i
    $quarterIndex = $i + 1
    $ws.Cells.Item($r, 1).Value = $oldA[$i].AddDays(1)
    $ws.Cells.Item($r, 2).Value = $newConsumption[$i]
    $ws.Cells.Item($r, 4).Value = "27.01.2026" + $quarterIndex
}
